$wb = $excel.ActiveWorkbook

# --- Sheet1: 展览 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 2451
$ws.Range("F3").Value = 719
$ws.Range("F4").Value = 239
$ws.Range("F5").Value = 410
$ws.Range("F6").Value = 683
$ws.Range("F8").Value = 886
$ws.Range("F9").Value = 563
$ws.Range("F10").Value = 927
$ws.Range("F13").Value = 436
$ws.Range("F16").Value = 1065
$ws.Range("F17").Value = 23884
$ws.Range("G17").Value = "已售罄"
$ws.Range("F18").Value = 2193
$ws.Range("F19").Value = 140
$ws.Range("F21").Value = 28
$ws.Range("F23").Value = 348
$ws.Range("F24").Value = 206
$ws.Range("F25").Value = 63
$ws.Range("F28").Value = 49
$ws.Range("F30").Value = 338
$ws.Range("F32").Value = 431
$ws.Range("F33").Value = 186
# --- Sheet2: 演出 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("F7").Value = 251
$ws.Range("F8").Value = 18
$ws.Range("F10").Value = 3593
$ws.Range("F12").Value = 143
$ws.Range("F16").Value = 15
$ws.Range("G18").Value = 380
$ws.Range("F19").Value = 4111
# --- Sheet3: 本地生活 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("F3").Value = 159
$ws.Range("F4").Value = 743
$ws.Range("F5").Value = 237
# --- Sheet4: 全部类型 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 159
$ws.Range("F4").Value = 2451
$ws.Range("F5").Value = 743
$ws.Range("F6").Value = 719
$ws.Range("F7").Value = 239
$ws.Range("F8").Value = 410
$ws.Range("F9").Value = 683
$ws.Range("F14").Value = 251
$ws.Range("F15").Value = 237
$ws.Range("F16").Value = 886
$ws.Range("F17").Value = 563
$ws.Range("F18").Value = 927
$ws.Range("F20").Value = 436
$ws.Range("F23").Value = 1065
$ws.Range("F24").Value = 23884
$ws.Range("G24").Value = "已售罄"
$ws.Range("F25").Value = 18
$ws.Range("F28").Value = 143
$ws.Range("F30").Value = 2193
$ws.Range("F31").Value = 140
$ws.Range("F33").Value = 28
$ws.Range("F36").Value = 348
$ws.Range("F37").Value = 206
$ws.Range("F38").Value = 63
$ws.Range("F41").Value = 15
$ws.Range("F42").Value = 49
$ws.Range("G45").Value = 380
$ws.Range("F47").Value = 431
$ws.Range("F48").Value = 186
